$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that used to sit
#    right after the title heading.
# ------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Giza Infinity Reels for Free
#    - Review & Demo") right before the last paragraph (the one that
#    used to contain the "Create a feature image..." image prompt),
#    and replace that last paragraph's text with the meta-description
#    copy, keeping its italic run formatting intact.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastRange = $d.Paragraphs.Item($count).Range
$lastRange.Collapse(1) | Out-Null

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Giza Infinity Reels for Free - Review &amp; Demo</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Explore the ancient Egyptian world of Giza Infinity Reels. Play for free and learn about symbol multipliers, free spins, and infinite ways to win.</w:t></w:r></w:p>'

$lastRange.InsertXML($xml) | Out-Null
